$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = [double]"0.9999999999999898"
$ws.Range("E2").Value = [double]"0.9999999999999898"

# Row 3
$ws.Range("D3").Value = [double]"7.353855186603038E-08"
$ws.Range("E3").Value = [double]"7.353855186603038E-08"

# Row 4
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = [double]"7.760516636441626E-09"
$ws.Range("E4").Value = [double]"7.760516636441626E-09"

# Row 5
$ws.Range("D5").Value = [double]"1.870927493784412E-06"
$ws.Range("E5").Value = [double]"1.870927493784412E-06"

# Row 6
$ws.Range("D6").Value = [double]"5.199034782839191E-38"
$ws.Range("E6").Value = [double]"5.199034782839191E-38"

# Row 7
$ws.Range("D7").Value = [double]"0.999999999995238"
$ws.Range("E7").Value = [double]"4.761968597222221E-12"

# Row 8
$ws.Range("D8").Value = [double]"0.9999999999322955"
$ws.Range("E8").Value = [double]"6.7704508666111E-11"

# Row 9
$ws.Range("D9").Value = [double]"0.5894363875520097"
$ws.Range("E9").Value = [double]"0.4105636124479903"

# Row 11
$ws.Range("D11").Value = [double]"0.0002267856319048671"
$ws.Range("E11").Value = [double]"0.9997732143680952"
$ws.Range("F11").Value = [double]"4.114037990570068"
$ws.Range("G11").Value = [double]"0.8"
